# Updated cryptos list (prices and 1h volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to remain
# Text (matching the source data which stores prices as literal strings),
# otherwise Excel would auto-convert them to numeric cells on assignment.
$textCells = $excel.Union($ws.Range("D5"), $ws.Range("D6"), $ws.Range("D7"), $ws.Range("D8"), $ws.Range("D16"), $ws.Range("D19"), $ws.Range("D20"), $ws.Range("D21"), $ws.Range("D22"), $ws.Range("D23"), $ws.Range("D24"), $ws.Range("D28"), $ws.Range("D30"), $ws.Range("D32"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D44"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D49"), $ws.Range("D51"))
$textCells.NumberFormat = "@"

$ws.Range("D2").Value = "61.816.55"
$ws.Range("E2").Value = "  -4.09%  "
$ws.Range("D3").Value = "2.987.19"
$ws.Range("E3").Value = "  -4.81%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "541.65"
$ws.Range("E5").Value = "  -5.55%  "
$ws.Range("D6").Value = "152.27"
$ws.Range("E6").Value = "  -7.34%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "2.998.05"
$ws.Range("E9").Value = "  -4.93%  "
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("E11").Value = "  -7.14%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").Value = "3.504.67"
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "61.833.43"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").Value = "23.96"
$ws.Range("E16").Value = "  -4.06%  "
$ws.Range("D17").Value = "2.986.96"
$ws.Range("E17").Value = "  -5.25%  "
$ws.Range("E18").Value = "  -5.45%  "
$ws.Range("D19").Value = "5.17"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "12.06"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").Value = "381.11"
$ws.Range("E21").Value = "  -6.53%  "
$ws.Range("D22").Value = "6.72"
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "66.01"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "3.106.71"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "0.0₃0943"
$ws.Range("E29").Value = "  -7.20%  "
$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  -8.40%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "20.52"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  -5.27%  "
$ws.Range("D34").Value = "159.72"
$ws.Range("E34").Value = "  -2.17%  "
$ws.Range("D35").Value = "5.93"
$ws.Range("E35").Value = "  -5.62%  "
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").Value = "  -5.91%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  -6.14%  "
$ws.Range("E39").Value = "  -8.29%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "2.419.78"
$ws.Range("E41").Value = "  -8.07%  "
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "22.03"
$ws.Range("E43").Value = "  -6.98%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.674"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").Value = "5.18"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "0.996"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "0.0245"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").Value = "19.84"
$ws.Range("E49").Value = "  -7.08%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "266.44"
$ws.Range("E51").Value = "  -8.33%  "
